$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows above row 881 to make room for a new weekly price block
# (this pushes the existing rows 881:958 down to 887:964, carrying their
# formatting - including the date style on column D - along with them).
$ws.Rows("881:886").Insert()

# Fixed/constant values shared by every data row in this block.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$catId     = 100112002
$categoria = "Pimiento"
$unidad    = "$/caja 18 kilos"
$origen    = "Provincia de Limarí"
$kgUnid    = 18
$clasif    = "Hortaliza"

$fecha = 44769

# New block of six rows (Primera/Segunda/Tercera for each variety) for the
# new reporting date, filling the rows just opened up above row 881.
$newRows = @(
    @{ Row=881; Variedad="Cuatro cascos verde"; Calidad="Primera"; Vol=800;  Min=25000; Max=26000; Prom=25500; PKg=1417 },
    @{ Row=882; Variedad="Cuatro cascos verde"; Calidad="Segunda"; Vol=900;  Min=22000; Max=23000; Prom=22500; PKg=1250 },
    @{ Row=883; Variedad="Cuatro cascos verde"; Calidad="Tercera"; Vol=700;  Min=19000; Max=20000; Prom=19500; PKg=1083 },
    @{ Row=884; Variedad="Morrón rojo";         Calidad="Primera"; Vol=1100; Min=28000; Max=29000; Prom=28500; PKg=1583 },
    @{ Row=885; Variedad="Morrón rojo";         Calidad="Segunda"; Vol=700;  Min=25000; Max=26000; Prom=25500; PKg=1417 },
    @{ Row=886; Variedad="Morrón rojo";         Calidad="Tercera"; Vol=700;  Min=22000; Max=23000; Prom=22500; PKg=1250 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $catId
    $ws.Cells.Item($row, 7).Value2  = $categoria
    $ws.Cells.Item($row, 8).Value2  = $r.Variedad
    $ws.Cells.Item($row, 9).Value2  = $r.Calidad
    $ws.Cells.Item($row, 10).Value2 = $r.Vol
    $ws.Cells.Item($row, 11).Value2 = $r.Min
    $ws.Cells.Item($row, 12).Value2 = $r.Max
    $ws.Cells.Item($row, 13).Value2 = $r.Prom
    $ws.Cells.Item($row, 14).Value2 = $unidad
    $ws.Cells.Item($row, 15).Value2 = $origen
    $ws.Cells.Item($row, 16).Value2 = $r.PKg
    $ws.Cells.Item($row, 17).Value2 = $kgUnid
    $ws.Cells.Item($row, 18).Value2 = $clasif
}
